$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2311.2222
$ws.Range("I106").Value = 2161.6924
$ws.Range("J106").Value = 2700
$ws.Range("K106").Value = 2161.6924
$ws.Range("L106").Value = 2700
$ws.Range("M106").Value = -1530.6924
$ws.Range("N106").Value = -3962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 709.85
$ws.Range("I74").Value = 715
$ws.Range("J74").Value = 700.2857
$ws.Range("K74").Value = 715
$ws.Range("L74").Value = 700.2857
$ws.Range("M74").Value = 159
$ws.Range("N74").Value = -2448.2857
$ws.Range("H77").Value = 709.85
$ws.Range("I77").Value = 715
$ws.Range("J77").Value = 700.2857
$ws.Range("K77").Value = 3575
$ws.Range("L77").Value = 3501.4285
$ws.Range("M77").Value = 793
$ws.Range("N77").Value = -12237.4285
$ws.Range("H97").Value = 1653.591
$ws.Range("I97").Value = 1670.091
$ws.Range("J97").Value = 1637.091
$ws.Range("K97").Value = 1670.091
$ws.Range("L97").Value = 1637.091
$ws.Range("M97").Value = -1174.091
$ws.Range("N97").Value = -2629.091
$ws.Range("H122").Value = 31078.115
$ws.Range("I122").Value = 2647.5173
$ws.Range("J122").Value = 168492.67
$ws.Range("K122").Value = 7942.5519
$ws.Range("L122").Value = 505478.01
$ws.Range("M122").Value = -5492.5519
$ws.Range("N122").Value = -510378.01

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 750
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = -360
$ws.Range("N8").Value = -1280
$ws.Range("H36").Value = 927.75
$ws.Range("I36").Value = 927.75
$ws.Range("K36").Value = 927.75
$ws.Range("M36").Value = -393.75
$ws.Range("H94").Value = 710.7857
$ws.Range("I94").Value = 743.375
$ws.Range("J94").Value = 667.3333
$ws.Range("K94").Value = 743.375
$ws.Range("L94").Value = 667.3333
$ws.Range("M94").Value = -292.375
$ws.Range("N94").Value = -1569.3333
$ws.Range("H99").Value = 1329.8667
$ws.Range("J99").Value = 2084.8333
$ws.Range("L99").Value = 2084.8333
$ws.Range("N99").Value = -5080.8333
$ws.Range("H112").Value = 10000
$ws.Range("I112").Value = 10000
$ws.Range("K112").Value = 10000
$ws.Range("M112").Value = -8523
$ws.Range("H134").Value = 21927.334
$ws.Range("I134").Value = 6818.25
$ws.Range("J134").Value = 142800
$ws.Range("K134").Value = 20454.75
$ws.Range("L134").Value = 428400
$ws.Range("M134").Value = -17919.75
$ws.Range("N134").Value = -433470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2445.5962
$ws.Range("I31").Value = 1710.6
$ws.Range("J31").Value = 3958.8235
$ws.Range("K31").Value = 1710.6
$ws.Range("L31").Value = 3958.8235
$ws.Range("M31").Value = -1415.6
$ws.Range("N31").Value = -4548.8235
$ws.Range("H34").Value = 2445.5962
$ws.Range("I34").Value = 1710.6
$ws.Range("J34").Value = 3958.8235
$ws.Range("K34").Value = 1710.6
$ws.Range("L34").Value = 3958.8235
$ws.Range("M34").Value = -1508.6
$ws.Range("N34").Value = -4362.8235
$ws.Range("H58").Value = 1031.0358
$ws.Range("I58").Value = 707.26086
$ws.Range("J58").Value = 2520.4
$ws.Range("K58").Value = 707.26086
$ws.Range("L58").Value = 2520.4
$ws.Range("M58").Value = -504.26086
$ws.Range("N58").Value = -2926.4
$ws.Range("H136").Value = 1031.0358
$ws.Range("I136").Value = 707.26086
$ws.Range("J136").Value = 2520.4
$ws.Range("K136").Value = 2121.78258
$ws.Range("L136").Value = 7561.200000000001
$ws.Range("M136").Value = 428.2174199999999
$ws.Range("N136").Value = -12661.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 830.25
$ws.Range("I51").Value = 144.125
$ws.Range("J51").Value = 2202.5
$ws.Range("K51").Value = 432.375
$ws.Range("L51").Value = 6607.5
$ws.Range("M51").Value = 27.625
$ws.Range("N51").Value = -7527.5
$ws.Range("H131").Value = 8197689
$ws.Range("J131").Value = 10417697
$ws.Range("L131").Value = 31253091
$ws.Range("N131").Value = -31263171
$ws.Range("H138").Value = 6367.727
$ws.Range("I138").Value = 6927.222
$ws.Range("J138").Value = 3850
$ws.Range("K138").Value = 20781.666
$ws.Range("L138").Value = 11550
$ws.Range("M138").Value = -15641.666
$ws.Range("N138").Value = -21830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2281.889
$ws.Range("I97").Value = 1940.6
$ws.Range("J97").Value = 2708.5
$ws.Range("K97").Value = 1940.6
$ws.Range("L97").Value = 2708.5
$ws.Range("M97").Value = -1444.6
$ws.Range("N97").Value = -3700.5
$ws.Range("H122").Value = 1870.3125
$ws.Range("I122").Value = 1904.1111
$ws.Range("J122").Value = 1826.8572
$ws.Range("K122").Value = 5712.3333
$ws.Range("L122").Value = 5480.571599999999
$ws.Range("M122").Value = -3262.3333
$ws.Range("N122").Value = -10380.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386
$ws.Range("H46").Value = 13650.25
$ws.Range("I46").Value = 933.3333
$ws.Range("J46").Value = 21280.4
$ws.Range("K46").Value = 933.3333
$ws.Range("L46").Value = 933.3333
$ws.Range("M46").Value = -745.3333
$ws.Range("N46").Value = -21656.4
$ws.Range("H80").Value = 28000
$ws.Range("J80").Value = 28000
$ws.Range("L80").Value = 28000
$ws.Range("N80").Value = -30246
$ws.Range("H83").Value = 28000
$ws.Range("J83").Value = 28000
$ws.Range("L83").Value = 84000
$ws.Range("N83").Value = -95232
$ws.Range("H100").Value = 2136
$ws.Range("I100").Value = 1465.7142
$ws.Range("J100").Value = 3178.6667
$ws.Range("K100").Value = 1465.7142
$ws.Range("L100").Value = 3178.6667
$ws.Range("M100").Value = -924.7141999999999
$ws.Range("N100").Value = -4260.6667
$ws.Range("H122").Value = 2154.8845
$ws.Range("I122").Value = 2194.7856
$ws.Range("J122").Value = 2108.3333
$ws.Range("K122").Value = 6584.3568
$ws.Range("L122").Value = 6324.999899999999
$ws.Range("M122").Value = -4134.3568
$ws.Range("N122").Value = -11224.9999
$ws.Range("H136").Value = 10515.765
$ws.Range("I136").Value = 12874.223
$ws.Range("K136").Value = 38622.669
$ws.Range("M136").Value = -36072.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 127250
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 203000
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 203000
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -205746
$ws.Range("H132").Value = 34662836
$ws.Range("I132").Value = 59530170
$ws.Range("J132").Value = 2024460.1
$ws.Range("K132").Value = 178590510
$ws.Range("L132").Value = 6073380.300000001
$ws.Range("M132").Value = -178587980
$ws.Range("N132").Value = -6078440.300000001
$ws.Range("H136").Value = 72639.57000000001
$ws.Range("I136").Value = 125821.75
$ws.Range("J136").Value = 1730
$ws.Range("K136").Value = 377465.25
$ws.Range("L136").Value = 5190
$ws.Range("N136").Value = -10290
